# 自动更新Excel文件
# For each data row (2..99), decrement the "剩余" (remaining, column E) by 1.
# If that would bring it down to 0, instead reset it to 10 and push the
# "开始时间" (start date, column F) forward by 10 (e.g. 20260114 -> 20260124).
# Row 36 already had remaining == 10 with a non-standard/odd start date and
# is intentionally left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F

    $eVal = $eCell.Value2
    if ($eVal -eq $null) {
        continue
    }

    if ($eVal -eq 10) {
        # Already at the "full / just refilled" value - leave untouched
        continue
    }

    if ($eVal -eq 1) {
        # Rolls over: remaining resets to 10, start date shifts 10 days later
        $eCell.Value2 = 10
        $fVal = $fCell.Value2
        if ($fVal -ne $null) {
            $fCell.Value2 = $fVal + 10
        }
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
